# Update NATMI LR-pair (Nppa-Npr3) sheet with new TPM-derived statistics.
# - Remove the three "Resolving-Mac" sending-cluster rows (old rows 11-13):
#   that cluster is no longer present in the refreshed TPM run.
# - Refresh the per-row ligand/receptor/edge expression statistics
#   (columns G-T) for the remaining 9 data rows with the new TPM values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old "Resolving-Mac" rows (11, 12, 13) entirely.
$ws.Range("A11:A13").EntireRow.Delete()

# New per-row values (columns G through T, i.e. columns 7-20) computed from
# the updated TPM data. Row keys below are the (sending, target) clusters,
# kept only as a readability comment - the sheet itself is unchanged in A-F.
$newValues = @{
    2  = @(0.1004046666666667, 0.301214, 0.3697907746891698, 0.44646652694238, 1, 0.3333333333333333, 0.01941266666666666, 0.058238, 0.01959774779686022, 0.01978237955823044, 0.001949122325777777, 0.017542100932, 0.00724706633996391, 0.008832170296019078)   # ECs -> ECs
    3  = @(0.1004046666666667, 0.301214, 0.3697907746891698, 0.44646652694238, 3, 1, 0.9434083333333333, 2.830225, 0.9524028256184742, 0.9613754796729409, 0.09472259923888887, 0.8525033931499999, 0.3521897787016099, 0.4292219714971426)                        # ECs -> FAPs
    4  = @(0.1004046666666667, 0.301214, 0.3697907746891698, 0.44646652694238, 1, 0.5, 0.027735, 0.05547, 0.02799942658466558, 0.01884214076882864, 0.00278472343, 0.01670834058, 0.01035392964759602, 0.008412385149218349)                                      # ECs -> MuSCs
    5  = @(0.03122233333333333, 0.093667, 0.1149919741207596, 0.1388354464902425, 1, 0.3333333333333333, 0.01941266666666666, 0.058238, 0.01959774779686022, 0.01978237955823044, 0.0006061087495555555, 0.005454978746, 0.002253583707481723, 0.002746495498606369) # FAPs -> ECs
    6  = @(0.03122233333333333, 0.093667, 0.1149919741207596, 0.1388354464902425, 3, 1, 0.9434083333333333, 2.830225, 0.9524028256184742, 0.9613754796729409, 0.02945540945277778, 0.265098685075, 0.1095186810760579, 0.1334729939651638)                        # FAPs -> FAPs
    7  = @(0.03122233333333333, 0.093667, 0.1149919741207596, 0.1388354464902425, 1, 0.5, 0.027735, 0.05547, 0.02799942658466558, 0.01884214076882864, 0.000865951415, 0.00519570849, 0.003219709337219973, 0.002615957026472326)                                  # FAPs -> MuSCs
    8  = @(0.1398905, 0.279781, 0.5152172511900707, 0.4146980265673775, 1, 0.3333333333333333, 0.01941266666666666, 0.058238, 0.01959774779686022, 0.01978237955823044, 0.002715647646333333, 0.016293885878, 0.01009709774941458, 0.008203713763604991)           # MuSCs -> ECs
    9  = @(0.1398905, 0.279781, 0.5152172511900707, 0.4146980265673775, 3, 1, 0.9434083333333333, 2.830225, 0.9524028256184742, 0.9613754796729409, 0.1319738634541667, 0.791843180725, 0.4906943658408065, 0.3986805142106345)                                    # MuSCs -> FAPs
    10 = @(0.1398905, 0.279781, 0.5152172511900707, 0.4146980265673775, 1, 0.5, 0.027735, 0.05547, 0.02799942658466558, 0.01884214076882864, 0.0038798630175, 0.01551945207, 0.01442578759984959, 0.007813798593137964)                                           # MuSCs -> MuSCs
}

foreach ($r in $newValues.Keys) {
    $rowVals = $newValues[$r]
    for ($i = 0; $i -lt $rowVals.Length; $i++) {
        $col = 7 + $i   # column G is index 7
        $ws.Cells.Item($r, $col).Value = $rowVals[$i]
    }
}
